$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the NaN student IDs in column A - renumber them
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15

# Add new "Projects Final Score" column header, bold
$ws.Range("C1").Value = "Projects Final Score"
$ws.Range("C1").Font.Bold = $true

# Update the selected cell, as saved by Excel
$ws.Range("H16").Select()
